$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: set the descriptive "Nombre"/"Descripcion" values (A and P columns) ---
$ws.Range("A2").Value = "N305-DELINE "
$ws.Range("P2").Value = "N305-DELINE "

$ws.Range("A3").Value = "N1035-ASDFA "
$ws.Range("P3").Value = "N1035-ASDFA "

$ws.Range("A4").Value = "N1036-N1036 -n1027"
$ws.Range("P4").Value = "N1036-N1036 -n1027"

$ws.Range("A5").Value = "N1037-ORGANIZADOR DE PLATO PARA LAVATORIO"
$ws.Range("P5").Value = "N1037-ORGANIZADOR DE PLATO PARA LAVATORIO"

$ws.Range("A6").Value = "N1038-TETERA DE ALUMINIO"
$ws.Range("P6").Value = "N1038-TETERA DE ALUMINIO"

# --- Step 2: set the short "Codigo Interno"/"Cod barras" values (B and T columns) ---
$ws.Range("B2").Value = "N305"
$ws.Range("T2").Value = "N305"

$ws.Range("B3").Value = "N1035"
$ws.Range("T3").Value = "N1035"

$ws.Range("B4").Value = "N1036"
$ws.Range("T4").Value = "N1036"

$ws.Range("B5").Value = "N1037"
$ws.Range("T5").Value = "N1037"

$ws.Range("B6").Value = "N1038"
$ws.Range("T6").Value = "N1038"

# --- Step 3: update numeric "Precio Unitario Venta" values ---
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 67

# --- Step 4: fill the rest of the new rows (4-6) with the remaining columns ---
$ws.Range("E4").Value = "NIU"
$ws.Range("F4").Value = "PEN"
$ws.Range("G4").Value = 8
$ws.Range("H4").Value = 10
$ws.Range("I4").Value = "SI"
$ws.Range("K4").Value = 10
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0

$ws.Range("E5").Value = "NIU"
$ws.Range("F5").Value = "PEN"
$ws.Range("G5").Value = 8
$ws.Range("H5").Value = 10
$ws.Range("I5").Value = "SI"
$ws.Range("K5").Value = 10
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0

$ws.Range("E6").Value = "NIU"
$ws.Range("F6").Value = "PEN"
$ws.Range("G6").Value = 8
$ws.Range("H6").Value = 10
$ws.Range("I6").Value = "SI"
$ws.Range("K6").Value = 10
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0

# --- Step 5: leave the active selection on A5, matching the saved view state ---
$ws.Range("A5").Select() | Out-Null
